$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G cells that contained the misspelled "impossibe" are corrected to
# "impossible" (an already-existing shared string). Once no cell references
# the misspelled string any more it is dropped from the shared strings table
# on save, which is exactly what the target workbook does.
$fixedCells = @(
    "G4","G6","G7","G8","G9","G10","G11","G12","G13","G16","G20",
    "G32","G34","G40","G41","G44","G45","G46","G47","G48","G49",
    "G51","G52","G53","G54","G55","G70","G71","G77","G78","G83",
    "G84","G89","G90","G91","G92","G96","G98","G99","G100","G101",
    "G103","G104","G105","G106","G107","G108","G109","G110","G112",
    "G113","G114","G117","G118","G119","G120"
)

foreach ($addr in $fixedCells) {
    $ws.Range($addr).Value = "impossible"
}

# Restore the view/selection that was active when the author saved the file.
$ws.Range("I46").Select()
